$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(9, 6).Value = 9814
$ws1.Cells.Item(11, 6).Value = 2658
$ws1.Cells.Item(13, 6).Value = 2393
$ws1.Cells.Item(16, 6).Value = 280
$ws1.Cells.Item(23, 6).Value = 302
$ws1.Cells.Item(32, 6).Value = 1685
$ws1.Cells.Item(33, 6).Value = 2828

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(14, 6).Value = 158

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(3, 6).Value = 953
$ws3.Cells.Item(5, 6).Value = 1768

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(4, 6).Value = 953
$ws4.Cells.Item(12, 6).Value = 9814
$ws4.Cells.Item(16, 6).Value = 2658
$ws4.Cells.Item(18, 6).Value = 2393
$ws4.Cells.Item(20, 6).Value = 280
$ws4.Cells.Item(26, 6).Value = 302
$ws4.Cells.Item(35, 6).Value = 1685
$ws4.Cells.Item(37, 6).Value = 2828
$ws4.Cells.Item(49, 6).Value = 158
